$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45180) that needs to
# be bumped to 45181 for every data row (rows 2 through 46).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
